$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.789.72"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.465.50"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.95"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.80%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  -1.43%  "

$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("E10").Value = "  -0.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.28"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.14"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000177"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.82%  "

$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.641.42"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.466.22"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.90"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.92"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "325.96"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.15%  "

$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.91%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +13.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.47"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "642.35"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.588.18"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0977"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.55%  "

$ws.Range("E29").Value = "  -13.24%  "

$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.94"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.13%  "

$ws.Range("E32").Value = "  -2.27%  "

$ws.Range("E33").Value = "  -4.08%  "

$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.368"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.46%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "151.08"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.60"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.32"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.70%  "

$ws.Range("E41").Value = "  -0.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.74"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₆0311"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -11.45%  "

$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "152.99"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.31"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.58"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.50%  "

$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.607"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0509"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.60%  "

$ws.Range("E51").Value = "  -1.28%  "
